# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) holds values like "5-12-2012-13" which should
# really read "2013-05-12" (ISO yyyy-mm-dd), stored as literal text.
#
# Excel auto-recognises strings that look like dates and silently turns
# them into date serial numbers as soon as they're typed into a cell, so
# we momentarily force the cell's number format to Text ("@") before
# writing the value, then clear the formatting back off again so the
# cell is left exactly as it started (no explicit style), with only its
# text content updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2013-05-12"
    $cell.ClearFormats()
}
